# chore: update Sheets via scheduled runner
# Refresh the cached market-board pricing columns (H, I/J, K/L) and the
# derived profit columns (M, N) on each job sheet (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) with newly scraped values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 83341016
$ws.Range("I64").Value = 190479460
$ws.Range("J64").Value = 11111.333
$ws.Range("K64").Value = 190479460
$ws.Range("L64").Value = 11111.333
$ws.Range("M64").Value = -190479212
$ws.Range("N64").Value = -11607.333
$ws.Range("H67").Value = 83341016
$ws.Range("I67").Value = 190479460
$ws.Range("J67").Value = 11111.333
$ws.Range("K67").Value = 190479460
$ws.Range("L67").Value = 11111.333
$ws.Range("M67").Value = -190478602
$ws.Range("N67").Value = -12827.333
$ws.Range("H74").Value = 14292732
$ws.Range("I74").Value = 23812920
$ws.Range("J74").Value = 12450
$ws.Range("K74").Value = 23812920
$ws.Range("L74").Value = 12450
$ws.Range("M74").Value = -23811984
$ws.Range("N74").Value = -14322
$ws.Range("H77").Value = 14292732
$ws.Range("I77").Value = 23812920
$ws.Range("J77").Value = 12450
$ws.Range("K77").Value = 119064600
$ws.Range("L77").Value = 62250
$ws.Range("M77").Value = -119059920
$ws.Range("N77").Value = -71610
$ws.Range("H107").Value = 1620.56
$ws.Range("I107").Value = 1061.7059
$ws.Range("J107").Value = 2808.125
$ws.Range("K107").Value = 1061.7059
$ws.Range("L107").Value = 2808.125
$ws.Range("M107").Value = 858.2941000000001
$ws.Range("N107").Value = -6648.125
$ws.Range("H112").Value = 3610.5588
$ws.Range("J112").Value = 3144.2122
$ws.Range("L112").Value = 9432.6366
$ws.Range("N112").Value = -11648.6366
$ws.Range("H127").Value = 5762
$ws.Range("I127").Value = 2990.1
$ws.Range("J127").Value = 8841.888999999999
$ws.Range("K127").Value = 8970.299999999999
$ws.Range("L127").Value = 26525.667
$ws.Range("M127").Value = -4010.299999999999
$ws.Range("N127").Value = -36445.667
$ws.Range("H135").Value = 5943.875
$ws.Range("I135").Value = 3170.0435
$ws.Range("K135").Value = 28530.3915
$ws.Range("M135").Value = -25995.3915
$ws.Range("H137").Value = 2365.054
$ws.Range("I137").Value = 837.5454999999999
$ws.Range("J137").Value = 4605.4
$ws.Range("K137").Value = 2512.6365
$ws.Range("L137").Value = 13816.2
$ws.Range("M137").Value = 37.36350000000039
$ws.Range("N137").Value = -18916.2
$ws.Range("H138").Value = 4955.6665
$ws.Range("J138").Value = 5769.4688
$ws.Range("L138").Value = 17308.4064
$ws.Range("N138").Value = -27588.4064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17931376
$ws.Range("I32").Value = 21283376
$ws.Range("J32").Value = 7428437.5
$ws.Range("K32").Value = 21283376
$ws.Range("L32").Value = 7428437.5
$ws.Range("M32").Value = -21283089
$ws.Range("N32").Value = -7429011.5
$ws.Range("H63").Value = 4601.3
$ws.Range("J63").Value = 5268
$ws.Range("L63").Value = 5268
$ws.Range("N63").Value = -6640
$ws.Range("H66").Value = 4601.3
$ws.Range("J66").Value = 5268
$ws.Range("L66").Value = 26340
$ws.Range("N66").Value = -33204
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""
$ws.Range("H139").Value = 97308.336
$ws.Range("J139").Value = 97308.336
$ws.Range("L139").Value = 97308.336
$ws.Range("N139").Value = -107588.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19330.285
$ws.Range("J82").Value = 23456
$ws.Range("L82").Value = 23456
$ws.Range("N82").Value = -24222
$ws.Range("H85").Value = 19330.285
$ws.Range("J85").Value = 23456
$ws.Range("L85").Value = 23456
$ws.Range("N85").Value = -26108
$ws.Range("H94").Value = 4266.921
$ws.Range("I94").Value = 2665.3928
$ws.Range("J94").Value = 8751.200000000001
$ws.Range("K94").Value = 2665.3928
$ws.Range("L94").Value = 8751.200000000001
$ws.Range("M94").Value = -2214.3928
$ws.Range("N94").Value = -9653.200000000001
$ws.Range("H99").Value = 8013.552
$ws.Range("I99").Value = 7260.2114
$ws.Range("K99").Value = 7260.2114
$ws.Range("M99").Value = -5762.2114
$ws.Range("H134").Value = 1407913.2
$ws.Range("I134").Value = 2034836.2
$ws.Range("J134").Value = 9392.538
$ws.Range("K134").Value = 6104508.6
$ws.Range("L134").Value = 28177.614
$ws.Range("M134").Value = -6101973.6
$ws.Range("N134").Value = -33247.614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 52635020
$ws.Range("I16").Value = 71429960
$ws.Range("J16").Value = 9199.4
$ws.Range("K16").Value = 71429960
$ws.Range("L16").Value = 9199.4
$ws.Range("M16").Value = -71429673
$ws.Range("N16").Value = -9773.4
$ws.Range("H31").Value = 8729.813
$ws.Range("I31").Value = 12880.353
$ws.Range("J31").Value = 7049.8335
$ws.Range("K31").Value = 12880.353
$ws.Range("L31").Value = 7049.8335
$ws.Range("M31").Value = -12585.353
$ws.Range("N31").Value = -7639.8335
$ws.Range("H34").Value = 8729.813
$ws.Range("I34").Value = 12880.353
$ws.Range("J34").Value = 7049.8335
$ws.Range("K34").Value = 12880.353
$ws.Range("L34").Value = 7049.8335
$ws.Range("M34").Value = -12678.353
$ws.Range("N34").Value = -7453.8335
$ws.Range("H51").Value = 40969
$ws.Range("I51").Value = 25000
$ws.Range("J51").Value = 46292
$ws.Range("K51").Value = 25000
$ws.Range("L51").Value = 46292
$ws.Range("M51").Value = -24264
$ws.Range("N51").Value = -47764
$ws.Range("H61").Value = 40969
$ws.Range("I61").Value = 25000
$ws.Range("J61").Value = 46292
$ws.Range("K61").Value = 25000
$ws.Range("L61").Value = 46292
$ws.Range("M61").Value = -24652
$ws.Range("N61").Value = -46988
$ws.Range("H107").Value = 670.1111
$ws.Range("I107").Value = 671.8333
$ws.Range("K107").Value = 671.8333
$ws.Range("M107").Value = 1248.1667
$ws.Range("H113").Value = 52635020
$ws.Range("I113").Value = 71429960
$ws.Range("J113").Value = 9199.4
$ws.Range("K113").Value = 71429960
$ws.Range("L113").Value = 9199.4
$ws.Range("M113").Value = -71427790
$ws.Range("N113").Value = -13539.4
$ws.Range("H132").Value = 6389.1113
$ws.Range("I132").Value = 4478.1714
$ws.Range("J132").Value = 9909.263000000001
$ws.Range("K132").Value = 13434.5142
$ws.Range("L132").Value = 29727.789
$ws.Range("M132").Value = -10904.5142
$ws.Range("N132").Value = -34787.789
$ws.Range("H141").Value = 181361.66
$ws.Range("J141").Value = 181361.66
$ws.Range("L141").Value = 181361.66
$ws.Range("N141").Value = -191721.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1557.25
$ws.Range("J11").Value = 133
$ws.Range("L11").Value = 399
$ws.Range("N11").Value = -679
$ws.Range("H22").Value = 20475
$ws.Range("I22").Value = 950.5
$ws.Range("K22").Value = 2851.5
$ws.Range("M22").Value = -2682.5
$ws.Range("H27").Value = 20475
$ws.Range("I27").Value = 950.5
$ws.Range("K27").Value = 2851.5
$ws.Range("M27").Value = -2749.5
$ws.Range("H43").Value = 9000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = ""
$ws.Range("H56").Value = 7788.2
$ws.Range("I56").Value = 7788.2
$ws.Range("K56").Value = 7788.2
$ws.Range("M56").Value = -7258.2
$ws.Range("H80").Value = 5175.8335
$ws.Range("I80").Value = 6500
$ws.Range("J80").Value = 4911
$ws.Range("K80").Value = 19500
$ws.Range("L80").Value = 14733
$ws.Range("M80").Value = -18564
$ws.Range("N80").Value = -16605
$ws.Range("H83").Value = 5175.8335
$ws.Range("I83").Value = 6500
$ws.Range("J83").Value = 4911
$ws.Range("K83").Value = 58500
$ws.Range("L83").Value = 44199
$ws.Range("M83").Value = -53820
$ws.Range("N83").Value = -53559
$ws.Range("H113").Value = 1530.6666
$ws.Range("I113").Value = 931.6667
$ws.Range("K113").Value = 2795.0001
$ws.Range("M113").Value = -625.0001000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6109.15
$ws.Range("J102").Value = 7432.5
$ws.Range("L102").Value = 7432.5
$ws.Range("N102").Value = -10676.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 20010
$ws.Range("J5").Value = 20010
$ws.Range("L5").Value = 20010
$ws.Range("N5").Value = -20236
$ws.Range("H46").Value = 35715784
$ws.Range("J46").Value = 62501800
$ws.Range("L46").Value = 62501800
$ws.Range("N46").Value = -62502176
$ws.Range("H68").Value = 2559.5
$ws.Range("I68").Value = 2079.3333
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2079.3333
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1330.3333
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 2559.5
$ws.Range("I71").Value = 2079.3333
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 10396.6665
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -6652.666499999999
$ws.Range("N71").Value = -27488
$ws.Range("H136").Value = 76939670
$ws.Range("I136").Value = 17011
$ws.Range("J136").Value = 111127530
$ws.Range("K136").Value = 51033
$ws.Range("L136").Value = 333382590
$ws.Range("M136").Value = -48483
$ws.Range("N136").Value = -333387690

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6795.71
$ws.Range("I132").Value = 5404.3447
$ws.Range("K132").Value = 16213.0341
$ws.Range("M132").Value = -13683.0341
$ws.Range("H136").Value = 13169158
$ws.Range("I136").Value = 17251340
$ws.Range("J136").Value = 15459.333
$ws.Range("K136").Value = 51754020
$ws.Range("L136").Value = 46377.999
$ws.Range("M136").Value = -51751470
$ws.Range("N136").Value = -51477.999
